# Weekly update: insert a new "Ajo" price record at the top of the data
# block (row 304), pushing the existing rows 304-332 down to 305-333.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 304 - this shifts rows 304:332 down to 305:333
# (and with them, the sheet's used-range dimension grows to A1:R333).
$ws.Rows.Item(304).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A304").Value = 5
$ws.Range("B304").Value = "Macroferia Regional de Talca"
$ws.Range("C304").Value = "Maule"
$ws.Range("D304").Value = 44769
$ws.Range("E304").Value = 7
$ws.Range("F304").Value = 100112003
$ws.Range("G304").Value = "Ajo"
$ws.Range("H304").Value = "Chino"
$ws.Range("I304").Value = "Primera"
$ws.Range("J304").Value = 300
$ws.Range("K304").Value = 28000
$ws.Range("L304").Value = 28000
$ws.Range("M304").Value = 28000
$ws.Range("N304").Value = "$/caja 10 kilos"
$ws.Range("O304").Value = "China"
$ws.Range("P304").Value = 2800
$ws.Range("Q304").Value = 10
$ws.Range("R304").Value = "Hortaliza"
